$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct tiny floating point drift on the existing row 9 timestamp
$ws.Range("A9").Value = 45809.39167575232

# Append new row 10 with the latest scraped price entry
$ws.Range("A10").Value = 45810.39396889304
$ws.Range("A10").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("B10").Value = "CREATINA MONOHIDRATO EN POLVO"
$ws.Range("C10").Value = "1Kg"
$ws.Range("D10").Value = "15,41€"
